# The title cell (A1) changes from "社員マスタ" to "社員一覧", and the
# subtitle cell (A2, "2024年度版") is removed, leaving row 2 blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "社員一覧"
$ws.Range("A2").ClearContents()

# Touch row 2 so Excel keeps an explicit (empty) <row r="2"/> element
# instead of dropping the row entirely when it has no cell content.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(2).Hidden = $false
